$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Weekly Coliflor price list update: a new sample (row 71) is inserted at the
# top of the data block; every existing row from 71..137 shifts its date /
# quality / volume / price data down by one row, and the row that used to be
# at 137 becomes a brand-new row 138 at the bottom.
# ---------------------------------------------------------------------------


# Row 71
$ws.Range("D71").Value = 44907
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 400
$ws.Range("L71").Value = 500
$ws.Range("M71").Value = 450
$ws.Range("P71").Value = 450

# Row 72
$ws.Range("D72").Value = 44411

# Row 73
$ws.Range("D73").Value = 44413
$ws.Range("I73").Value = "Tercera"
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 500
$ws.Range("L73").Value = 600
$ws.Range("M73").Value = 550
$ws.Range("P73").Value = 550

# Row 74
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 700
$ws.Range("K74").Value = 800
$ws.Range("L74").Value = 900
$ws.Range("M74").Value = 850
$ws.Range("P74").Value = 850

# Row 75
$ws.Range("D75").Value = 44350
$ws.Range("I75").Value = "Tercera"
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 500
$ws.Range("L75").Value = 600
$ws.Range("M75").Value = 550
$ws.Range("P75").Value = 550

# Row 76
$ws.Range("I76").Value = "Segunda"
$ws.Range("J76").Value = 900
$ws.Range("K76").Value = 800
$ws.Range("L76").Value = 1000
$ws.Range("M76").Value = 900
$ws.Range("P76").Value = 900

# Row 77
$ws.Range("D77").Value = 44406
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 500
$ws.Range("L77").Value = 600
$ws.Range("M77").Value = 550
$ws.Range("P77").Value = 550

# Row 78
$ws.Range("D78").Value = 44873
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 600
$ws.Range("L78").Value = 700
$ws.Range("M78").Value = 650
$ws.Range("P78").Value = 650

# Row 79
$ws.Range("D79").Value = 44837
$ws.Range("I79").Value = "Tercera"
$ws.Range("K79").Value = 700
$ws.Range("L79").Value = 750
$ws.Range("M79").Value = 725
$ws.Range("P79").Value = 725

# Row 80
$ws.Range("I80").Value = "Segunda"
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 600
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = 733
$ws.Range("P80").Value = 733

# Row 81
$ws.Range("D81").Value = 44348
$ws.Range("I81").Value = "Tercera"
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 700
$ws.Range("L81").Value = 700
$ws.Range("M81").Value = 700
$ws.Range("P81").Value = 700

# Row 82
$ws.Range("I82").Value = "Primera"
$ws.Range("K82").Value = 950
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = 975
$ws.Range("P82").Value = 975

# Row 83
$ws.Range("D83").Value = 44274
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 700
$ws.Range("L83").Value = 800
$ws.Range("M83").Value = 750
$ws.Range("P83").Value = 750

# Row 84
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 700
$ws.Range("K84").Value = 800
$ws.Range("L84").Value = 900
$ws.Range("M84").Value = 850
$ws.Range("P84").Value = 850

# Row 85
$ws.Range("D85").Value = 44427
$ws.Range("I85").Value = "Tercera"
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 700
$ws.Range("M85").Value = 650
$ws.Range("P85").Value = 650

# Row 86
$ws.Range("I86").Value = "Segunda"
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = 950
$ws.Range("P86").Value = 950

# Row 87
$ws.Range("D87").Value = 44334
$ws.Range("J87").Value = 800
$ws.Range("K87").Value = 500
$ws.Range("L87").Value = 600
$ws.Range("M87").Value = 550
$ws.Range("P87").Value = 550

# Row 88
$ws.Range("D88").Value = 44894
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 350
$ws.Range("L88").Value = 450
$ws.Range("M88").Value = 400
$ws.Range("P88").Value = 400

# Row 89
$ws.Range("D89").Value = 44383
$ws.Range("I89").Value = "Tercera"
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 500
$ws.Range("L89").Value = 600
$ws.Range("M89").Value = 540
$ws.Range("P89").Value = 540

# Row 90
$ws.Range("I90").Value = "Segunda"
$ws.Range("J90").Value = 800
$ws.Range("K90").Value = 700
$ws.Range("L90").Value = 800
$ws.Range("M90").Value = 750
$ws.Range("P90").Value = 750

# Row 91
$ws.Range("D91").Value = 44882
$ws.Range("I91").Value = "Tercera"
$ws.Range("J91").Value = 700
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 600
$ws.Range("M91").Value = 550
$ws.Range("P91").Value = 550

# Row 92
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 850
$ws.Range("L92").Value = 950
$ws.Range("M92").Value = 900
$ws.Range("P92").Value = 900

# Row 93
$ws.Range("I93").Value = "Segunda"
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 650
$ws.Range("L93").Value = 750
$ws.Range("M93").Value = 700
$ws.Range("P93").Value = 700

# Row 94
$ws.Range("D94").Value = 44285
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 400
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = 450
$ws.Range("P94").Value = 450

# Row 95
$ws.Range("D95").Value = 44426
$ws.Range("I95").Value = "Tercera"
$ws.Range("J95").Value = 700
$ws.Range("K95").Value = 900
$ws.Range("L95").Value = 1000
$ws.Range("M95").Value = 950
$ws.Range("P95").Value = 950

# Row 96
$ws.Range("I96").Value = "Segunda"
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = 1100
$ws.Range("P96").Value = 1100

# Row 97
$ws.Range("D97").Value = 44390
$ws.Range("I97").Value = "Tercera"
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = 650
$ws.Range("P97").Value = 650

# Row 98
$ws.Range("I98").Value = "Segunda"
$ws.Range("J98").Value = 1300
$ws.Range("K98").Value = 700
$ws.Range("L98").Value = 750
$ws.Range("M98").Value = 725
$ws.Range("P98").Value = 725

# Row 99
$ws.Range("D99").Value = 44159
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 650
$ws.Range("M99").Value = 625
$ws.Range("P99").Value = 625

# Row 100
$ws.Range("D100").Value = 44567
$ws.Range("I100").Value = "Tercera"
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 400
$ws.Range("L100").Value = 500
$ws.Range("M100").Value = 450
$ws.Range("P100").Value = 450

# Row 101
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 1000
$ws.Range("M101").Value = 950
$ws.Range("P101").Value = 950

# Row 102
$ws.Range("D102").Value = 44496
$ws.Range("J102").Value = 700
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = 750
$ws.Range("P102").Value = 750

# Row 103
$ws.Range("I103").Value = "Segunda"
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 700
$ws.Range("L103").Value = 750
$ws.Range("M103").Value = 725
$ws.Range("P103").Value = 725

# Row 104
$ws.Range("D104").Value = 44273
$ws.Range("J104").Value = 900

# Row 105
$ws.Range("D105").Value = 44385
$ws.Range("J105").Value = 700
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 600
$ws.Range("M105").Value = 550
$ws.Range("P105").Value = 550

# Row 106
$ws.Range("D106").Value = 44608
$ws.Range("J106").Value = 800
$ws.Range("K106").Value = 1000
$ws.Range("L106").Value = 1200
$ws.Range("M106").Value = 1100
$ws.Range("P106").Value = 1100

# Row 107
$ws.Range("D107").Value = 44327
$ws.Range("I107").Value = "Tercera"
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 550
$ws.Range("P107").Value = 550

# Row 108
$ws.Range("I108").Value = "Segunda"
$ws.Range("J108").Value = 500
$ws.Range("K108").Value = 600
$ws.Range("L108").Value = 700
$ws.Range("M108").Value = 650
$ws.Range("P108").Value = 650

# Row 109
$ws.Range("D109").Value = 44474
$ws.Range("J109").Value = 1200
$ws.Range("K109").Value = 400
$ws.Range("M109").Value = 450
$ws.Range("P109").Value = 450

# Row 110
$ws.Range("D110").Value = 44377
$ws.Range("I110").Value = "Tercera"
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 450
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 480
$ws.Range("P110").Value = 480

# Row 111
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 1200
$ws.Range("K111").Value = 500
$ws.Range("L111").Value = 600
$ws.Range("M111").Value = 550
$ws.Range("P111").Value = 550

# Row 112
$ws.Range("D112").Value = 44558
$ws.Range("I112").Value = "Tercera"
$ws.Range("K112").Value = 350
$ws.Range("L112").Value = 400
$ws.Range("M112").Value = 375
$ws.Range("P112").Value = 375

# Row 113
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 650
$ws.Range("P113").Value = 650

# Row 114
$ws.Range("D114").Value = 44523
$ws.Range("K114").Value = 400
$ws.Range("L114").Value = 500
$ws.Range("M114").Value = 450
$ws.Range("P114").Value = 450

# Row 115
$ws.Range("I115").Value = "Segunda"
$ws.Range("K115").Value = 600
$ws.Range("L115").Value = 700
$ws.Range("M115").Value = 650
$ws.Range("P115").Value = 650

# Row 116
$ws.Range("D116").Value = 44280
$ws.Range("I116").Value = "Tercera"
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 500
$ws.Range("L116").Value = 600
$ws.Range("M116").Value = 550
$ws.Range("P116").Value = 550

# Row 117
$ws.Range("I117").Value = "Segunda"
$ws.Range("J117").Value = 1000
$ws.Range("K117").Value = 650
$ws.Range("L117").Value = 750
$ws.Range("M117").Value = 700
$ws.Range("P117").Value = 700

# Row 118
$ws.Range("D118").Value = 44642
$ws.Range("I118").Value = "Tercera"
$ws.Range("J118").Value = 1200
$ws.Range("K118").Value = 500
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 550
$ws.Range("P118").Value = 550

# Row 119
$ws.Range("I119").Value = "Segunda"
$ws.Range("J119").Value = 700
$ws.Range("K119").Value = 800
$ws.Range("L119").Value = 900
$ws.Range("M119").Value = 850
$ws.Range("P119").Value = 850

# Row 120
$ws.Range("D120").Value = 44432
$ws.Range("I120").Value = "Tercera"
$ws.Range("J120").Value = 500

# Row 121
$ws.Range("I121").Value = "Segunda"
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 700
$ws.Range("M121").Value = 650
$ws.Range("P121").Value = 650

# Row 122
$ws.Range("D122").Value = 44294
$ws.Range("I122").Value = "Tercera"
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 450
$ws.Range("L122").Value = 500
$ws.Range("M122").Value = 475
$ws.Range("P122").Value = 475

# Row 123
$ws.Range("I123").Value = "Segunda"
$ws.Range("J123").Value = 500
$ws.Range("K123").Value = 800
$ws.Range("L123").Value = 900
$ws.Range("M123").Value = 850
$ws.Range("P123").Value = 850

# Row 124
$ws.Range("D124").Value = 44803
$ws.Range("I124").Value = "Tercera"
$ws.Range("J124").Value = 800
$ws.Range("K124").Value = 500
$ws.Range("L124").Value = 600
$ws.Range("M124").Value = 550
$ws.Range("P124").Value = 550

# Row 125
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 1300
$ws.Range("K125").Value = 900
$ws.Range("L125").Value = 1000
$ws.Range("M125").Value = 946
$ws.Range("P125").Value = 946

# Row 126
$ws.Range("D126").Value = 44881
$ws.Range("J126").Value = 850
$ws.Range("M126").Value = 679
$ws.Range("P126").Value = 679

# Row 127
$ws.Range("I127").Value = "Segunda"
$ws.Range("J127").Value = 800
$ws.Range("K127").Value = 650
$ws.Range("L127").Value = 700
$ws.Range("M127").Value = 675
$ws.Range("P127").Value = 675

# Row 128
$ws.Range("D128").Value = 44323
$ws.Range("I128").Value = "Tercera"
$ws.Range("J128").Value = 1500
$ws.Range("K128").Value = 500
$ws.Range("L128").Value = 600
$ws.Range("M128").Value = 550
$ws.Range("P128").Value = 550

# Row 129
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 1200
$ws.Range("K129").Value = 750
$ws.Range("L129").Value = 800
$ws.Range("M129").Value = 775
$ws.Range("P129").Value = 775

# Row 130
$ws.Range("D130").Value = 44537
$ws.Range("I130").Value = "Tercera"
$ws.Range("J130").Value = 1000
$ws.Range("K130").Value = 600
$ws.Range("L130").Value = 650
$ws.Range("M130").Value = 625
$ws.Range("P130").Value = 625

# Row 131
$ws.Range("I131").Value = "Segunda"
$ws.Range("J131").Value = 900
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 1000
$ws.Range("M131").Value = 950
$ws.Range("P131").Value = 950

# Row 132
$ws.Range("D132").Value = 44644
$ws.Range("I132").Value = "Tercera"
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 500
$ws.Range("L132").Value = 600
$ws.Range("M132").Value = 550
$ws.Range("P132").Value = 550

# Row 133
$ws.Range("D133").Value = 44544
$ws.Range("J133").Value = 1200
$ws.Range("K133").Value = 400
$ws.Range("L133").Value = 500
$ws.Range("M133").Value = 450
$ws.Range("P133").Value = 450

# Row 134
$ws.Range("I134").Value = "Segunda"
$ws.Range("K134").Value = 700
$ws.Range("L134").Value = 800
$ws.Range("M134").Value = 750
$ws.Range("P134").Value = 750

# Row 135
$ws.Range("D135").Value = 44336
$ws.Range("I135").Value = "Tercera"
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 400
$ws.Range("L135").Value = 500
$ws.Range("M135").Value = 450
$ws.Range("P135").Value = 450

# Row 136
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 900
$ws.Range("L136").Value = 1000
$ws.Range("M136").Value = 950
$ws.Range("P136").Value = 950

# Row 137
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 900
$ws.Range("K137").Value = 600
$ws.Range("L137").Value = 700
$ws.Range("M137").Value = 650
$ws.Range("P137").Value = 650

# Row 138
$ws.Range("D138").Value = 44649
$ws.Range("D138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I138").Value = "Tercera"
$ws.Range("J138").Value = 1200
$ws.Range("K138").Value = 400
$ws.Range("L138").Value = 500
$ws.Range("M138").Value = 450
$ws.Range("P138").Value = 450

# Row 138 is a brand-new row -- fill in the columns that stay constant for
# every record in this subset (mercado/región/categoría/unidad/origen/etc.)
$ws.Range("A138").Value = 1
$ws.Range("B138").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C138").Value = "Arica y Parinacota"
$ws.Range("E138").Value = 15
$ws.Range("F138").Value = 100112008
$ws.Range("G138").Value = "Coliflor"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("N138").Value = "$/unidad"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = "Hortaliza"

